# Update "想去人数" (F) and "最低票价" (G) figures on both the "展览" and
# "全部类型" sheets to reflect newly scraped numbers (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 2645
    $ws.Range("F4").Value = 540
    $ws.Range("F6").Value = 6605
    $ws.Range("F7").Value = 457
    $ws.Range("F9").Value = 9

    $ws.Range("F10").Value = 42
    $ws.Range("G10").Value = 50
}

# "全部类型" sheet received one additional correction on row 8 that the
# "展览" sheet did not.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 9
